$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This paragraph (in the "Drawables" section) currently reads:
#   "...smallest-screen-width." <br/> "Note: Sun-tiles...1dp border" "?"
#   <bookmarkStart _GoBack/><bookmarkEnd/>
#
# Target layout:
#   "...smallest-screen-width." + NEW "Advantage of using..." run
#   <bookmarkStart _GoBack/><bookmarkEnd/>            (bookmark MOVED here)
#   <br/> "Note: Sun-tiles...1dp border" "?"
#   + NEW " Also if ..." + NEW " is set" + NEW ", the Ima" + NEW "geView..." runs
#
# NOTE on ordering: this runtime re-coalesces adjacent same-formatted runs
# across a paragraph whenever an edit is applied to that paragraph, *unless*
# a bookmark sits exactly at the join point (which acts as a hard barrier).
# So: do ALL text insertions first (while the pre-existing _GoBack bookmark
# is still parked at its original spot, protecting the "?"-adjacent new
# runs from being re-absorbed), and only relocate the bookmark as the very
# last step.
# ---------------------------------------------------------------------------

# 1) Add the new sentence right after "...smallest-screen-width."
$rngAdvantage = $d.Content
$rngAdvantage.Find.Execute("smallest-screen-width.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngAdvantage.Collapse(0)
$rngAdvantage.InsertAfter(" Advantage of using “wrap_content” is that the layout can remain the same while there are different drawables which there needs to be anyway.")

# 2) Add the four new trailing runs right after the existing "?" run, which
#    is still immediately followed by the (not-yet-moved) _GoBack bookmark.
$rngTail = $d.Content
$rngTail.Find.Execute("within the 1dp border?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngTail.Collapse(0)

$rngTail.InsertAfter(" Also if “wrap_content” is used and no image")
$rngTail.Collapse(0)

$rngTail.InsertAfter(" is set")
$rngTail.Collapse(0)

$rngTail.InsertAfter(", the Ima")
$rngTail.Collapse(0)

$rngTail.InsertAfter("geView is shrunk to a 0x0 or 1x1 box, but resizes correctly when image resource is set. May need to have ‘border’ drawable that is a PNG rather than a shape XML file.")

# 3) Finally, relocate the _GoBack bookmark from the end of the paragraph to
#    right after the newly inserted "Advantage..." sentence.
$d.Bookmarks.Item("_GoBack").Delete()

$rngBookmark = $d.Content
$rngBookmark.Find.Execute("there needs to be anyway.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngBookmark.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngBookmark)
